$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 4 - HuobiToken price update
Set-TextValue "D4" "5.418"

# Row 5 - Cronos price update
Set-TextValue "D5" "0.05960"

# Row 6 - GateToken price update
Set-TextValue "D6" "3.388"

# Row 7 - MXToken price update
Set-TextValue "D7" "0.8085"

# Row 8 - FTXToken price update
Set-TextValue "D8" "0.9217"

# Row 9 - now One (was WazirX)
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D9" "0.01117"
$ws.Range("E9").Value = "8OneONEBestin24h"

# Row 10 - now WazirX (was MandalaExchangeToken)
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1417"
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11 - now MandalaExchangeToken (was LiechtensteinCryptoassetsExchange)
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07428"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12 - now LiechtensteinCryptoassetsExchange (was BitrueCoin)
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03412"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13 - now BitrueCoin (was BitMartToken)
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03055"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14 - now BitMartToken (was MCDex)
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09348"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15 - now MCDex (was BitForexToken)
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.932"
$ws.Range("E15").Value = "14MCDexMCB"

# Row 16 - now BitForexToken (was CoinExToken)
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001595"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17 - now CoinExToken (was One)
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04810"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18 - TigerCash price update
Set-TextValue "D18" "0.005474"

# Row 19 - HotbitToken price update
Set-TextValue "D19" "0.004158"

# Row 20 - BitKan price update
Set-TextValue "D20" "0.0009810"

# Row 21 - NitroEx price update
Set-TextValue "D21" "0.00007703"

# Row 22 - LEO price update
Set-TextValue "D22" "3.660"

# Row 23 - KuCoinToken price update
Set-TextValue "D23" "6.433"

# Row 40 - IDEX price update
Set-TextValue "D40" "0.03933"

# Row 41 - now KickToken (was BKEXToken)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006193"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 - now BKEXToken (was CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1073"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 - now CEJI (was KickToken)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders price update and label change
Set-TextValue "D44" "0.007166"
$ws.Range("E44").Value = "43LocalTradersLCT"

# Row 45 - CoinLion price update
Set-TextValue "D45" "0.00005193"

# Row 47 - ACDXExchange price update
Set-TextValue "D47" "0.0005802"

# Row 48 - CoinbaseStockToken price update
Set-TextValue "D48" "1.050"
